$d = $word.ActiveDocument

function Find-Range([string]$text) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $text"
    }
    return $rng
}

# ---------------------------------------------------------------------------
# 1) "To some extent ... supportive culture. After all, we believe that our
#    success comes down to our diligence and determination." loses its final
#    sentence (and the lone "." run that follows it), leaving the paragraph
#    ending "...supportive culture. " (trailing space retained/added).
# ---------------------------------------------------------------------------
$runA = Find-Range("To some extent, it was surprising to find that we functioned very well as a team. According to our personality tests, we were a relatively introverted group of individuals. We had concerned that the fact would affect our collaboration adversely. We were proud of overcoming the concern and forming a well-organized team with a positive and supportive culture. After all, we believe that our success comes down to our diligence and determination")
$periodRun = $d.Range($runA.End, $runA.End + 1)
$periodRun.Text = ""
$runA.Text = "To some extent, it was surprising to find that we functioned very well as a team. According to our personality tests, we were a relatively introverted group of individuals. We had concerned that the fact would affect our collaboration adversely. We were proud of overcoming the concern and forming a well-organized team with a positive and supportive culture. "

# ---------------------------------------------------------------------------
# 2) "Each member surprised us ... Brandon has excellent " / "interview
#    skills, ... Hugo has amazing people skill..." -- the phrase "interview
#    skills, which surprised not only the team but also himself. Hugo has
#    amazing " moves from the start of the run following the page break to
#    the end of the run preceding it.
# ---------------------------------------------------------------------------
$runAfterBreak = Find-Range("interview skills, which surprised not only the team but also himself. Hugo has amazing people skill at such a young age. Taylen always surprises us with how skillful he is with IT.  Tim is a well-balanced businessperson who is logical, thoughtful, and assertive. Tetsu is a hard-working individual who still enjoys studying in his mid-40s.")
$runAfterBreak.Text = "people skill at such a young age. Taylen always surprises us with how skillful he is with IT.  Tim is a well-balanced businessperson who is logical, thoughtful, and assertive. Tetsu is a hard-working individual who still enjoys studying in his mid-40s."

$runBeforeBreak = Find-Range("Brandon has excellent ")
$runBeforeBreak.Text = "Brandon has excellent interview skills, which surprised not only the team but also himself. Hugo has amazing "

# ---------------------------------------------------------------------------
# 3) The closing paragraph's only sentence is removed entirely, leaving an
#    empty paragraph in its place.
# ---------------------------------------------------------------------------
$closing = Find-Range("We look forward to completing the assignment together and developing an even better team for the rest of the course.")
$closing.Text = ""
